$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Status" (column F) values for several defects, and highlight the
#     ones that just changed to Fixed/fixed with the existing orange fill that
#     the sheet already uses (the same fill previously used on I7's "fixed"
#     comment). Grab the orange-fill template from I7 before it gets
#     repurposed below. ---

# F7: status moves from "Open" to "fixed" (same text that used to live in the
# Comments cell I7); pick up I7's existing orange-fill/vertical-center format.
$ws.Range("I7").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = "fixed"

# F6: "Open" -> "Fixed", highlighted with the orange fill (no vertical
# centering on this row).
$ws.Range("I7").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").VerticalAlignment = -4107
$ws.Range("F6").Value = "Fixed"

# F18: "Open" -> "fixed", same orange highlight as F6.
$ws.Range("I7").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").VerticalAlignment = -4107
$ws.Range("F18").Value = "fixed"

# F22: "Open" -> "Fixed", same orange highlight.
$ws.Range("I7").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").VerticalAlignment = -4107
$ws.Range("F22").Value = "Fixed"

# F25: "Open" -> "Fixed", same orange highlight.
$ws.Range("I7").Copy()
$ws.Range("F25").PasteSpecial(-4122)
$ws.Range("F25").VerticalAlignment = -4107
$ws.Range("F25").Value = "Fixed"

# I7's old "fixed" comment is no longer needed now that the status itself says
# so - clear it, but keep the row's plain (non-filled) vertically centered
# look (same formatting already used elsewhere in row 7, e.g. H7). Do this
# last since the steps above use I7's current (orange) format as a template.
$ws.Range("H7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").ClearContents()

# F13 / F23: text-only status updates, formatting stays as-is.
$ws.Range("F13").Value = "fixed"
$ws.Range("F23").Value = "open"

$excel.CutCopyMode = 0

# --- Update the current selection/scroll position left behind in the sheet
#     view (cosmetic, matches where the author was last working) ---
$ws.Range("F15").Select()

Write-Output "done"
